# Update column F (dSF) values on Sheet1 to match the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -9
    3  = -4
    4  = -9
    5  = -5
    6  = -4
    7  = 4
    10 = -7
    13 = 2
    14 = -9
    16 = -5
    17 = -7
    18 = -3
    19 = -2
    22 = -2
    23 = 0
    24 = -4
    26 = 2
    27 = -1
    29 = 1
    30 = 0
    33 = 1
    38 = -3
    39 = 2
    42 = -5
    44 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
